$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().EndsWith("16")) {
        $cell.Value2 = $val.ToString().Substring(0, $val.ToString().Length - 2)
    }
}
